$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("事件功能")

# Mark the previously "not started" rows (失效控件, 生效控件, 偏移控件, 闪烁控件, 停止闪烁)
# as completed now that enable/disableElement, moveElement and blinkElement have landed.
$ws.Range("B13").Value = "已完成"
$ws.Range("B14").Value = "已完成"
$ws.Range("B15").Value = "已完成"
$ws.Range("B16").Value = "已完成"
$ws.Range("B17").Value = "已完成"

# Highlight the "偏移控件" (moveElement) row with the same top/left accent
# border + centred alignment used to call out freshly landed rows.
$ws.Range("A15").Borders.Item(7).LineStyle = 1
$ws.Range("A15").Borders.Item(7).ThemeColor = 7
$ws.Range("A15").Borders.Item(8).LineStyle = 1
$ws.Range("A15").Borders.Item(8).ThemeColor = 7
$ws.Range("A15").HorizontalAlignment = -4108
$ws.Range("A15").VerticalAlignment = -4108

$ws.Range("B17").Select()
